$d = $word.ActiveDocument

# =====================================================================
# Change 1: Insert a new list-item paragraph "Make sure they are all in
# the correct folders!" right before the "If there are hundreds of
# changed files..." paragraph (same list level as that paragraph).
# =====================================================================
$hundredsIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*hundreds of changed files*") {
        $hundredsIdx = $i
    }
}

$hundredsPara = $d.Paragraphs($hundredsIdx)
$hundredsPara.Range.InsertParagraphBefore()
$d.Paragraphs($hundredsIdx).Range.Text = "Make sure they are all in the correct folders!"

# =====================================================================
# Change 2: Append a new run to the "discard changes...classroom
# documents." paragraph, then move the "_GoBack" bookmark (currently at
# the very end of the document) to the end of this paragraph, right
# after the newly appended run.
# =====================================================================

# Remove the pre-existing "_GoBack" bookmark (it sits at the end of the
# document's last paragraph in the original file).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$discardIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*discard changes*") {
        $discardIdx = $i
    }
}

# Split point: just before the paragraph mark of the "discard changes"
# paragraph. Drop a throwaway bookmark there first so that the text
# inserted next becomes its own run (matching the original authoring)
# instead of being merged into the preceding run.
$discardPara = $d.Paragraphs($discardIdx)
$splitRange = $discardPara.Range.Duplicate
[void]$splitRange.MoveEnd(1, -1)
$splitRange.Collapse(0)
$d.Bookmarks.Add("TempRunSplit", $splitRange)

$splitRange.InsertAfter(" Do not delete or edit other student" + [char]8217 + "s work.")

if ($d.Bookmarks.Exists("TempRunSplit")) {
    $d.Bookmarks("TempRunSplit").Delete()
}

# Now place the real "_GoBack" bookmark at the very end of that same
# paragraph (after the run we just added, before the paragraph mark).
# A collapsed bookmark placed exactly one character before a paragraph
# mark gets mis-anchored to the whole paragraph, so temporarily type a
# placeholder character, anchor the bookmark just before it, and then
# delete the placeholder again.
$discardPara2 = $d.Paragraphs($discardIdx)
$endRange = $discardPara2.Range.Duplicate
[void]$endRange.MoveEnd(1, -1)
$endRange.Collapse(0)
$endRange.InsertAfter("X")

$discardPara3 = $d.Paragraphs($discardIdx)
$bookmarkRange = $discardPara3.Range.Duplicate
[void]$bookmarkRange.MoveEnd(1, -2)
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$discardPara4 = $d.Paragraphs($discardIdx)
$placeholderRange = $discardPara4.Range.Duplicate
$placeholderRange.Start = $placeholderRange.End - 2
$placeholderRange.End = $placeholderRange.End - 1
$placeholderRange.Delete()
